$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("G2").Value = 15.475219
$ws.Range("H2").Value = 46.425657
$ws.Range("I2").Value = 0.1968226394800303
$ws.Range("J2").Value = 0.1968226394800303
$ws.Range("M2").Value = 0.063209
$ws.Range("N2").Value = 0.189627
$ws.Range("O2").Value = 0.03474243690088823
$ws.Range("P2").Value = 0.03474243690088823
$ws.Range("Q2").Value = 0.9781731177710001
$ws.Range("R2").Value = 8.803558059939
$ws.Range("S2").Value = 0.006838098132801224
$ws.Range("T2").Value = 0.006838098132801225

$ws.Range("G3").Value = 15.475219
$ws.Range("H3").Value = 46.425657
$ws.Range("I3").Value = 0.1968226394800303
$ws.Range("J3").Value = 0.1968226394800303
$ws.Range("N3").Value = 4.918502999999999
$ws.Range("O3").Value = 0.9011416102365667
$ws.Range("P3").Value = 0.9011416102365668
$ws.Range("Q3").Value = 25.371637025719
$ws.Range("R3").Value = 228.344733231471
$ws.Range("S3").Value = 0.1773650702720457
$ws.Range("T3").Value = 0.1773650702720458

$ws.Range("G4").Value = 15.475219
$ws.Range("H4").Value = 46.425657
$ws.Range("I4").Value = 0.1968226394800303
$ws.Range("J4").Value = 0.1968226394800303
$ws.Range("O4").Value = 0.06411595286254507
$ws.Range("P4").Value = 0.06411595286254508
$ws.Range("Q4").Value = 1.80518429635
$ws.Range("R4").Value = 16.24665866715
$ws.Range("S4").Value = 0.01261947107518332
$ws.Range("T4").Value = 0.01261947107518333

$ws.Range("I5").Value = 0.6648050370145543
$ws.Range("J5").Value = 0.6648050370145544
$ws.Range("M5").Value = 0.063209
$ws.Range("N5").Value = 0.189627
$ws.Range("O5").Value = 0.03474243690088823
$ws.Range("P5").Value = 0.03474243690088823
$ws.Range("Q5").Value = 3.303961462382333
$ws.Range("R5").Value = 29.735653161441
$ws.Range("S5").Value = 0.02309694704987082
$ws.Range("T5").Value = 0.02309694704987082

$ws.Range("I6").Value = 0.6648050370145543
$ws.Range("J6").Value = 0.6648050370145544
$ws.Range("N6").Value = 4.918502999999999
$ws.Range("O6").Value = 0.9011416102365667
$ws.Range("P6").Value = 0.9011416102365668
$ws.Range("Q6").Value = 85.69741842992767
$ws.Range("R6").Value = 771.2767658693489
$ws.Range("S6").Value = 0.5990834815486757
$ws.Range("T6").Value = 0.599083481548676

$ws.Range("I7").Value = 0.6648050370145543
$ws.Range("J7").Value = 0.6648050370145544
$ws.Range("O7").Value = 0.06411595286254507
$ws.Range("P7").Value = 0.06411595286254508
$ws.Range("S7").Value = 0.04262460841600769
$ws.Range("T7").Value = 0.04262460841600771

$ws.Range("I8").Value = 0.1383723235054153
$ws.Range("J8").Value = 0.1383723235054153
$ws.Range("M8").Value = 0.063209
$ws.Range("N8").Value = 0.189627
$ws.Range("O8").Value = 0.03474243690088823
$ws.Range("P8").Value = 0.03474243690088823
$ws.Range("Q8").Value = 0.6876855602286667
$ws.Range("R8").Value = 6.189170042058
$ws.Range("S8").Value = 0.004807391718216185
$ws.Range("T8").Value = 0.004807391718216185

$ws.Range("I9").Value = 0.1383723235054153
$ws.Range("J9").Value = 0.1383723235054153
$ws.Range("N9").Value = 4.918502999999999
$ws.Range("O9").Value = 0.9011416102365667
$ws.Range("P9").Value = 0.9011416102365668
$ws.Range("S9").Value = 0.1246930584158451
$ws.Range("T9").Value = 0.1246930584158451

$ws.Range("I10").Value = 0.1383723235054153
$ws.Range("J10").Value = 0.1383723235054153
$ws.Range("O10").Value = 0.06411595286254507
$ws.Range("P10").Value = 0.06411595286254508
$ws.Range("S10").Value = 0.008871873371354045
$ws.Range("T10").Value = 0.008871873371354047
